# Update the division problems in the practice worksheet table.
# Each cell's "a÷b=" text is replaced with the new problem per the commit.

$d = $word.ActiveDocument

$d.Content.Find.Execute("26÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "16÷8=", 2) | Out-Null
$d.Content.Find.Execute("30÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "75÷8=", 2) | Out-Null
$d.Content.Find.Execute("16÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "41÷7=", 2) | Out-Null
$d.Content.Find.Execute("76÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "25÷2=", 2) | Out-Null
$d.Content.Find.Execute("42÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "34÷2=", 2) | Out-Null
$d.Content.Find.Execute("63÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "53÷9=", 2) | Out-Null
$d.Content.Find.Execute("83÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "67÷8=", 2) | Out-Null
$d.Content.Find.Execute("49÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "46÷7=", 2) | Out-Null
$d.Content.Find.Execute("41÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "47÷8=", 2) | Out-Null
$d.Content.Find.Execute("86÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "52÷6=", 2) | Out-Null
$d.Content.Find.Execute("84÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "80÷5=", 2) | Out-Null
$d.Content.Find.Execute("41÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "54÷4=", 2) | Out-Null
$d.Content.Find.Execute("77÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "57÷9=", 2) | Out-Null
$d.Content.Find.Execute("10÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "34÷9=", 2) | Out-Null
$d.Content.Find.Execute("86÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "70÷4=", 2) | Out-Null
$d.Content.Find.Execute("61÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "21÷8=", 2) | Out-Null
$d.Content.Find.Execute("43÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "44÷9=", 2) | Out-Null
$d.Content.Find.Execute("27÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "96÷7=", 2) | Out-Null
$d.Content.Find.Execute("69÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "55÷5=", 2) | Out-Null
$d.Content.Find.Execute("40÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "87÷3=", 2) | Out-Null
$d.Content.Find.Execute("55÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "19÷8=", 2) | Out-Null
$d.Content.Find.Execute("63÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "29÷5=", 2) | Out-Null
$d.Content.Find.Execute("20÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "26÷4=", 2) | Out-Null
$d.Content.Find.Execute("89÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "33÷8=", 2) | Out-Null
$d.Content.Find.Execute("42÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "52÷9=", 2) | Out-Null

Write-Output "Replacements complete."
